# Fruta / hortaliza, semanal
# Insert a new weekly price record at the top of the "Terminal La Palmera de
# La Serena - Cilantro" data block (row 214), pushing the existing rows
# 214:225 down to 215:226. The sheet's used range grows from A1:R225 to
# A1:R226 automatically as part of the row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 214; this shifts rows 214:225 to
# 215:226 and carries each cell's formatting (e.g. the date style on column D).
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with the new weekly record.
$ws.Cells.Item(214, 1).Value  = 8
$ws.Cells.Item(214, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(214, 3).Value  = "Coquimbo"
$ws.Cells.Item(214, 4).Value  = 45041
$ws.Cells.Item(214, 5).Value  = 4
$ws.Cells.Item(214, 6).Value  = 100112040
$ws.Cells.Item(214, 7).Value  = "Cilantro"
$ws.Cells.Item(214, 8).Value  = "Sin especificar"
$ws.Cells.Item(214, 9).Value  = "Primera"
$ws.Cells.Item(214, 10).Value = 2600
$ws.Cells.Item(214, 11).Value = 2000
$ws.Cells.Item(214, 12).Value = 2500
$ws.Cells.Item(214, 13).Value = 2250
$ws.Cells.Item(214, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(214, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(214, 16).Value = 1500
$ws.Cells.Item(214, 17).Value = 1.5
$ws.Cells.Item(214, 18).Value = "Hortaliza"
